$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = [double]"0.8469462330186518"
$ws.Range("E2").Value = [double]"0.8469462330186518"

# Row 3
$ws.Range("D3").Value = [double]"2.737654710891644E-06"
$ws.Range("E3").Value = [double]"2.737654710891644E-06"

# Row 4
$ws.Range("D4").Value = [double]"0.002307396152455066"
$ws.Range("E4").Value = [double]"0.002307396152455066"

# Row 5
$ws.Range("D5").Value = [double]"0.0004722175605528189"
$ws.Range("E5").Value = [double]"0.0004722175605528189"

# Row 6
$ws.Range("D6").Value = [double]"0.863860727118141"
$ws.Range("E6").Value = [double]"0.863860727118141"

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = [double]"2.281701924783946E-06"
$ws.Range("E7").Value = [double]"0.9999977182980753"

# Row 8
$ws.Range("D8").Value = [double]"0.978147963702534"
$ws.Range("E8").Value = [double]"0.02185203629746602"

# Row 9
$ws.Range("D9").Value = [double]"0.9846908339042469"
$ws.Range("E9").Value = [double]"0.01530916609575306"

# Row 10
$ws.Range("D10").Value = [double]"0.9999999999972069"
$ws.Range("E10").Value = [double]"2.793099085351969E-12"

# Row 11
$ws.Range("D11").Value = [double]"1"
$ws.Range("E11").Value = [double]"0"
$ws.Range("F11").Value = [double]"1.690193772315979"
$ws.Range("G11").Value = [double]"0.7"

# Row 12
$ws.Range("D12").Value = [double]"0.9359363654338638"
$ws.Range("E12").Value = [double]"0.9359363654338638"

# Row 13
$ws.Range("D13").Value = [double]"1.791891961738539E-08"
$ws.Range("E13").Value = [double]"1.791891961738539E-08"

# Row 14
$ws.Range("D14").Value = [double]"0.0004539061076841042"
$ws.Range("E14").Value = [double]"0.0004539061076841042"

# Row 15
$ws.Range("D15").Value = [double]"4.003747610730699E-05"
$ws.Range("E15").Value = [double]"4.003747610730699E-05"

# Row 16
$ws.Range("D16").Value = [double]"0.9175364866847586"
$ws.Range("E16").Value = [double]"0.9175364866847586"

# Row 17
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = [double]"1.037066595485808E-07"
$ws.Range("E17").Value = [double]"0.9999998962933404"

# Row 18
$ws.Range("D18").Value = [double]"0.9867543559950352"
$ws.Range("E18").Value = [double]"0.01324564400496475"

# Row 19
$ws.Range("D19").Value = [double]"0.9927786964501278"
$ws.Range("E19").Value = [double]"0.00722130354987216"

# Row 20
$ws.Range("D20").Value = [double]"1"
$ws.Range("E20").Value = [double]"0"

# Row 21
$ws.Range("D21").Value = [double]"1"
$ws.Range("E21").Value = [double]"0"
$ws.Range("F21").Value = [double]"2.134605646133423"
$ws.Range("G21").Value = [double]"0.7"
